# PALI rNPV Model - bullish FSCD cPoS projection update
# ------------------------------------------------------------------
# 1) FSCD Phase Transition Probabilities (rNPV Model sheet) revised upward
# 2) Sensitivity "D23"/"D24" FSCD cPoS scenario inputs revised upward to match
# 3) B102 formula on rNPV Model corrected to reference B79 (FSCD NDA->Approval probability)
#    instead of B80 (the $ NPV line) when summing strategic/optionality value per share
# 4) Text edits:
#    - "Diluted Shares"!A8 label shortened (warrants fully exercised note removed)
#    - "Sensitivity"!A32 competitor-risk note rewritten
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$wsDiluted = $wb.Worksheets.Item("Diluted Shares")
$wsModel   = $wb.Worksheets.Item("rNPV Model")
$wsSens    = $wb.Worksheets.Item("Sensitivity")

# --- 1) FSCD Phase Transition Probabilities (rNPV Model!B39:B42) ---
$wsModel.Range("B39").Value = 0.85   # Phase 1b -> Phase 2
$wsModel.Range("B40").Value = 0.35   # Phase 2 -> Phase 3
$wsModel.Range("B41").Value = 0.7    # Phase 3 -> NDA Filing
$wsModel.Range("B42").Value = 0.9    # NDA -> FDA Approval

# --- 2) Sensitivity scenario FSCD cPoS inputs (D23:D24) ---
$wsSens.Range("D23").Value = 0.315
$wsSens.Range("D24").Value = 0.315

# --- 3) Fix strategic premium per-share formula to use B79 (FSCD NDA probability)
#        instead of B80 ---
$wsModel.Range("B102").Formula = "=(B79+B81)*1000000/B91"

# --- 4) Text updates (order matters for shared-string table append order) ---
$wsSens.Range("A32").Value = "Competitor: ALK5 class history likely necessitates FDA-mandated cardiac monitoring (REMS) for Agomab despite clean Ph2. PALI-2108 (gut-restricted prodrug) offers unmonitored safety profile, capturing First-Line volume"
$wsDiluted.Range("A8").Value = "Oct 2025 Pre-Funded Warrants"

# --- 5) Scroll position bookkeeping (matches author's viewport when saving) ---
$wsDiluted.Application.ActiveWindow.ScrollRow = 4
$wsModel.Application.ActiveWindow.ScrollRow = 68

$wb.Application.Calculate()
